$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 mirrors the style/formatting of the other header cells (e.g. E1)
$ws.Range("F1").Value = "19-jun"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Fill F2:F25 with placeholder "-" values (no special style, like B2:E25)
for ($row = 2; $row -le 25; $row++) {
    $ws.Cells.Item($row, 6).Value = "-"
}
